$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing header/amount layout - scraper now only writes part number + description
$ws.Cells.Clear()

# Row 1: part number, bold, wrap text
$ws.Range("A1").Value = "3273114"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").WrapText = $true

$ws.Range("B1").Value = ""

# Row 2: description, bold, wrap text
$ws.Range("A2").Value = "6SL32105BE211UV0"
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").WrapText = $true

$ws.Range("A2").Select()
